$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1 (paragraph "{m:if self.name <>'anydsl'}"):
#   insert a new run containing a single space between the "<>" run
#   and the "'" run, so the text reads "... <> 'anydsl'}".
# -----------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2Start = $p2.Range.Start

# Locate "<>" inside this paragraph and insert a space right after it.
$find2 = $d.Range($p2Start, $p2.Range.End)
$find2.Find.Execute("<>", $true, $false, $false, $false, $false, $true, 1, $false, "<> ", 2)

# The paragraph text is now merged into a single run by the text edit
# above; re-split it right after the inserted space so the space lives
# in its own run (matching the target structure) without altering any
# visible formatting.
$p2 = $d.Paragraphs.Item(2)
$splitAt = $p2.Range.Start + 18
$spaceRange = $d.Range($splitAt, $splitAt + 1)
$spaceRange.Font.Bold = $true
$spaceRange.Font.Bold = $false

# -----------------------------------------------------------------
# Change 2 (paragraph "{m:elseif self.name = 'anydsl'}"):
#   split the run holding "elseif self.name = 'anydsl'}" into two
#   runs: "elseif self.name = 'anydsl'" and "}".
# -----------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$p4End = $p4.Range.End
$closeBrace = $d.Range($p4End - 2, $p4End - 1)
$closeBrace.Font.Bold = $true
$closeBrace.Font.Bold = $false
